$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font + border) from the existing H1 header cell
# onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill column I with 1s and column J as a copy of column H for data rows 2..30
for ($r = 2; $r -le 30; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
